$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3663.6667
$ws.Range("J17").Value = 3663.6667
$ws.Range("L17").Value = 10991.0001
$ws.Range("N17").Value = -11327.0001

$ws.Range("H28").Value = 1125.5
$ws.Range("I28").Value = 1125.5
$ws.Range("K28").Value = 1125.5
$ws.Range("M28").Value = -640.5

$ws.Range("H47").Value = 39999
$ws.Range("I47").Value = 39999
$ws.Range("K47").Value = 39999
$ws.Range("M47").Value = -39027

$ws.Range("H70").Value = 1461.3334
$ws.Range("I70").Value = 1749.8334
$ws.Range("J70").Value = 1172.8334
$ws.Range("K70").Value = 5249.5002
$ws.Range("L70").Value = 3518.5002
$ws.Range("M70").Value = -4979.5002
$ws.Range("N70").Value = -4058.5002

$ws.Range("H73").Value = 1461.3334
$ws.Range("I73").Value = 1749.8334
$ws.Range("J73").Value = 1172.8334
$ws.Range("K73").Value = 5249.5002
$ws.Range("L73").Value = 3518.5002
$ws.Range("M73").Value = -4313.5002
$ws.Range("N73").Value = -5390.5002

$ws.Range("H86").Value = 5024.375
$ws.Range("I86").Value = 2539.2
$ws.Range("J86").Value = 9166.333000000001
$ws.Range("K86").Value = 2539.2
$ws.Range("L86").Value = 9166.333000000001
$ws.Range("M86").Value = -1416.2
$ws.Range("N86").Value = -11412.333

$ws.Range("H89").Value = 5024.375
$ws.Range("I89").Value = 2539.2
$ws.Range("J89").Value = 9166.333000000001
$ws.Range("K89").Value = 12696
$ws.Range("L89").Value = 45831.665
$ws.Range("M89").Value = -7080
$ws.Range("N89").Value = -57063.665

$ws.Range("H92").Value = 779.1429000000001
$ws.Range("J92").Value = 671
$ws.Range("L92").Value = 671
$ws.Range("N92").Value = -3167

$ws.Range("H98").Value = 3473
$ws.Range("I98").Value = 3591.25
$ws.Range("K98").Value = 3591.25
$ws.Range("M98").Value = -2093.25

$ws.Range("H113").Value = 3174.6924
$ws.Range("I113").Value = 2378
$ws.Range("K113").Value = 2378
$ws.Range("M113").Value = 876

$ws.Range("H122").Value = 3473
$ws.Range("I122").Value = 3591.25
$ws.Range("K122").Value = 10773.75
$ws.Range("M122").Value = -8323.75

$ws.Range("H129").Value = 2059.1226
$ws.Range("J129").Value = 2076.9368
$ws.Range("L129").Value = 6230.8104
$ws.Range("N129").Value = -16230.8104

$ws.Range("H132").Value = 1624.8914
$ws.Range("I132").Value = 1255.909
$ws.Range("K132").Value = 3767.727
$ws.Range("M132").Value = -1237.727

$ws.Range("H138").Value = 2610.1667
$ws.Range("I138").Value = 1455.6
$ws.Range("K138").Value = 4366.799999999999
$ws.Range("M138").Value = 773.2000000000007

$ws.Range("H141").Value = 55080.527
$ws.Range("I141").Value = 64490.688
$ws.Range("K141").Value = 193472.064
$ws.Range("M141").Value = -188292.064

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3237.0488
$ws.Range("I32").Value = 3302.975
$ws.Range("J32").Value = 600
$ws.Range("K32").Value = 3302.975
$ws.Range("L32").Value = 600
$ws.Range("M32").Value = -3015.975
$ws.Range("N32").Value = -1174

$ws.Range("H45").Value = 3689.6667
$ws.Range("I45").Value = 1538.8
$ws.Range("J45").Value = 5226
$ws.Range("K45").Value = 1538.8
$ws.Range("L45").Value = 5226
$ws.Range("M45").Value = -1161.8
$ws.Range("N45").Value = -5980

$ws.Range("H60").Value = 16299.5
$ws.Range("I60").Value = 16299.5
$ws.Range("K60").Value = 16299.5
$ws.Range("M60").Value = -15566.5

$ws.Range("H74").Value = 1920.3055
$ws.Range("J74").Value = 1741.3889
$ws.Range("L74").Value = 1741.3889
$ws.Range("N74").Value = -3489.3889

$ws.Range("H77").Value = 1920.3055
$ws.Range("J77").Value = 1741.3889
$ws.Range("L77").Value = 8706.9445
$ws.Range("N77").Value = -17442.9445

$ws.Range("H122").Value = 3111
$ws.Range("I122").Value = 3016.25
$ws.Range("K122").Value = 9048.75
$ws.Range("M122").Value = -6598.75

$ws.Range("H132").Value = 2036.3636
$ws.Range("I132").Value = 1708.5
$ws.Range("K132").Value = 5125.5
$ws.Range("M132").Value = -2595.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 121523
$ws.Range("I9").Value = 200051
$ws.Range("J9").Value = 42995
$ws.Range("K9").Value = 200051
$ws.Range("L9").Value = 42995
$ws.Range("M9").Value = -199883
$ws.Range("N9").Value = -43331

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H86").Value = 7068.4614
$ws.Range("I86").Value = 2117.4
$ws.Range("J86").Value = 10162.875
$ws.Range("K86").Value = 2117.4
$ws.Range("L86").Value = 10162.875
$ws.Range("M86").Value = -994.4000000000001
$ws.Range("N86").Value = -12408.875

$ws.Range("H89").Value = 7068.4614
$ws.Range("I89").Value = 2117.4
$ws.Range("J89").Value = 10162.875
$ws.Range("K89").Value = 10587
$ws.Range("L89").Value = 50814.375
$ws.Range("M89").Value = -4971
$ws.Range("N89").Value = -62046.375

$ws.Range("H94").Value = 4496.6665
$ws.Range("I94").Value = 6434
$ws.Range("K94").Value = 6434
$ws.Range("M94").Value = -5983

$ws.Range("H107").Value = 2603
$ws.Range("I107").Value = 2954.75
$ws.Range("K107").Value = 2954.75
$ws.Range("M107").Value = -1034.75

$ws.Range("H134").Value = 2277.1304
$ws.Range("I134").Value = 2192.628
$ws.Range("K134").Value = 6577.884
$ws.Range("M134").Value = -4042.884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1199.5
$ws.Range("I122").Value = 1199.5
$ws.Range("K122").Value = 3598.5
$ws.Range("M122").Value = -1148.5

$ws.Range("H134").Value = 2873.5
$ws.Range("I134").Value = 2327.457
$ws.Range("J134").Value = 4610.909
$ws.Range("K134").Value = 6982.370999999999
$ws.Range("L134").Value = 13832.727
$ws.Range("M134").Value = -4447.370999999999
$ws.Range("N134").Value = -18902.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H111").Value = 526
$ws.Range("I111").Value = 526
$ws.Range("K111").Value = 1578
$ws.Range("M111").Value = 1489

$ws.Range("H116").Value = 127851.695
$ws.Range("I116").Value = 154057.3
$ws.Range("K116").Value = 462171.9
$ws.Range("M116").Value = -458729.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 100000
$ws.Range("J94").Value = 100000
$ws.Range("L94").Value = 100000
$ws.Range("N94").Value = -101352

$ws.Range("H99").Value = 7233
$ws.Range("I99").Value = 7233
$ws.Range("K99").Value = 7233
$ws.Range("M99").Value = -4987

$ws.Range("H107").Value = 319.6
$ws.Range("I107").Value = 32.666668
$ws.Range("K107").Value = 32.666668
$ws.Range("M107").Value = 1887.333332

$ws.Range("H113").Value = 2006.625
$ws.Range("I113").Value = 1675.5
$ws.Range("K113").Value = 1675.5
$ws.Range("M113").Value = 494.5

$ws.Range("H122").Value = 3812.7222
$ws.Range("I122").Value = 4856.222
$ws.Range("J122").Value = 2769.2222
$ws.Range("K122").Value = 14568.666
$ws.Range("L122").Value = 8307.6666
$ws.Range("M122").Value = -12118.666
$ws.Range("N122").Value = -13207.6666

$ws.Range("H126").Value = 9335.727999999999
$ws.Range("I126").Value = 3670.4285
$ws.Range("J126").Value = 19250
$ws.Range("K126").Value = 11011.2855
$ws.Range("L126").Value = 57750
$ws.Range("M126").Value = -8541.2855
$ws.Range("N126").Value = -62690

$ws.Range("H132").Value = 3192.0557
$ws.Range("I132").Value = 3210.4546
$ws.Range("J132").Value = 3163.1428
$ws.Range("K132").Value = 9631.363799999999
$ws.Range("L132").Value = 9489.428400000001
$ws.Range("M132").Value = -7101.363799999999
$ws.Range("N132").Value = -14549.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2535.8333
$ws.Range("I7").Value = 2304
$ws.Range("K7").Value = 2304
$ws.Range("M7").Value = -2192

$ws.Range("H93").Value = 52627.25
$ws.Range("I93").Value = 2131.75
$ws.Range("J93").Value = 77875
$ws.Range("K93").Value = 2131.75
$ws.Range("L93").Value = 77875
$ws.Range("M93").Value = -883.75
$ws.Range("N93").Value = -80371

$ws.Range("H126").Value = 2535.8333
$ws.Range("I126").Value = 2304
$ws.Range("K126").Value = 6912
$ws.Range("M126").Value = -4442

$ws.Range("H132").Value = 15194
$ws.Range("I132").Value = 11996.333
$ws.Range("J132").Value = 19990.5
$ws.Range("K132").Value = 35988.999
$ws.Range("L132").Value = 59971.5
$ws.Range("M132").Value = -33458.999
$ws.Range("N132").Value = -65031.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1568.5
$ws.Range("I107").Value = 1200.25
$ws.Range("K107").Value = 3600.75
$ws.Range("M107").Value = -1680.75

$ws.Range("H122").Value = 9977.333000000001
$ws.Range("I122").Value = 9977.333000000001
$ws.Range("K122").Value = 29931.999
$ws.Range("M122").Value = -27481.999

$ws.Range("H126").Value = 2799.375
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
